$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "Koti"
$ws.Range("B2").Value = "images/kotesh.png"

# Row 3 updates
$ws.Range("A3").Value = "Vaggg"
$ws.Range("B3").Value = "images/Vagdevi(photo).png"

# Remove row 4 entirely
$ws.Range("A4:B4").Delete()
